$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 18-hole course data to append: Masterton Golf Course
$courseName = "Masterton Golf Course"
$holes = @(
    @(1, 5, 410),
    @(2, 5, 420),
    @(3, 3, 187),
    @(4, 4, 225),
    @(5, 4, 326),
    @(6, 4, 360),
    @(7, 4, 406),
    @(8, 3, 114),
    @(9, 4, 356),
    @(10, 4, 357),
    @(11, 3, 159),
    @(12, 4, 298),
    @(13, 3, 129),
    @(14, 4, 305),
    @(15, 4, 318),
    @(16, 4, 326),
    @(17, 4, 327),
    @(18, 5, 453)
)

$startRow = 29
$row = $startRow
foreach ($hole in $holes) {
    $ws.Cells.Item($row, 1).Value = $courseName
    $ws.Cells.Item($row, 2).Value = $hole[0]
    $ws.Cells.Item($row, 3).Value = $hole[1]
    $ws.Cells.Item($row, 4).Value = $hole[2]
    $row = $row + 1
}

# Resize column A to fit the newly-added, longer course name
$ws.Columns.Item(1).AutoFit()

# Update the view: scroll down and select the cell the author landed on
$null = $ws.Range("I24").Select()
